# Horarios actualizados - Linea 141 - 2026-01-11 - 12:11:52 scrape run
# Applies the new scrape results (rows re-fetched + new arrivals appended)
# to all three sheets: LP1912, LP1912-215, 6203-6173.

$wb = $excel.ActiveWorkbook

function Set-RowData($ws, $rows) {
    foreach ($r in $rows) {
        $rn = $r[0]
        if ($r[1] -ne $null) { $ws.Cells.Item($rn, 1).Value = $r[1] }
        if ($r[2] -ne $null) { $ws.Cells.Item($rn, 2).Value = $r[2] }
        if ($r[3] -ne $null) { $ws.Cells.Item($rn, 3).Value = $r[3] }
        if ($r[4] -ne $null) { $ws.Cells.Item($rn, 4).Value = $r[4] }
        if ($r[5] -ne $null) { $ws.Cells.Item($rn, 5).Value = $r[5] }
    }
}

# --- Sheet "LP1912" ---
# Columns: Hora_Scrap | Hora_Llegada | Linea | Minutos | Parada
$ws1 = $wb.Worksheets.Item("LP1912")
$rows1 = @(
    ,@(2, "Última actualización: 12:11:52", $null, $null, $null, $null)
    ,@(3, "Total filas: 157", $null, $null, $null, $null)
    ,@(23, "05:57:13", "07:21", "23_HERNANDEZ", 84, "LP1912")
    ,@(24, "06:17:28", "07:21", "16_SANTA ANA", 64, "LP1912")
    ,@(33, "06:17:28", "08:00", "17_ROMERO", 103, "LP1912")
    ,@(34, "06:46:50", "08:00", "16_SANTA ANA", 74, "LP1912")
    ,@(64, "07:50:33", "09:39", "15_ABASTO", 109, "LP1912")
    ,@(65, "08:27:16", "09:39", "23_HERNANDEZ", 72, "LP1912")
    ,@(122, "12:11:52", "12:11", "16_P MOR-SANTA ANA", 0, "LP1912")
    ,@(123, "12:11:52", "12:13", "15_ABASTO", 2, "LP1912")
    ,@(124, "10:37:52", "12:16", "10_OLMOS", 99, "LP1912")
    ,@(125, "11:11:33", "12:17", "10_OLMOS", 66, "LP1912")
    ,@(126, "10:37:52", "12:21", "215C_EL PATO", 104, "LP1912")
    ,@(127, "11:11:33", "12:22", "215C_EL PATO", 71, "LP1912")
    ,@(128, "11:47:17", "12:32", "23_HERNANDEZ", 45, "LP1912")
    ,@(129, "10:37:52", "12:32", "14_ABASTO", 115, "LP1912")
    ,@(130, "11:34:59", "12:33", "15_ABASTO", 59, "LP1912")
    ,@(131, "11:47:17", "12:33", "14_ABASTO", 46, "LP1912")
    ,@(132, "10:37:52", "12:34", "15_ABASTO", 117, "LP1912")
    ,@(133, "11:11:33", "12:35", "23_HERNANDEZ", 84, "LP1912")
    ,@(134, "11:34:59", "12:35", "27_EL RETIRO", 61, "LP1912")
    ,@(135, "10:50:41", "12:36", "27_EL RETIRO", 106, "LP1912")
    ,@(136, "11:34:59", "12:36", "23_HERNANDEZ", 62, "LP1912")
    ,@(137, "11:47:17", "12:37", "27_EL RETIRO", 50, "LP1912")
    ,@(138, "11:52:01", "12:37", "23_HERNANDEZ", 45, "LP1912")
    ,@(139, "11:34:59", "12:47", "16_SANTA ANA", 73, "LP1912")
    ,@(140, "11:34:59", "12:47", "14_ABASTO", 73, "LP1912")
    ,@(141, "11:34:59", "12:47", "15X38_ABASTO", 73, "LP1912")
    ,@(142, "11:47:17", "12:48", "14_ABASTO", 61, "LP1912")
    ,@(143, "11:11:33", "12:48", "15X38_ABASTO", 97, "LP1912")
    ,@(144, "10:50:41", "12:48", "16_SANTA ANA", 118, "LP1912")
    ,@(145, "11:11:33", "13:02", "11_ETCHEVERRY", 111, "LP1912")
    ,@(146, "11:34:59", "13:03", "215C_EL PATO", 89, "LP1912")
    ,@(147, "11:47:17", "13:03", "11_ETCHEVERRY", 76, "LP1912")
    ,@(148, "11:47:17", "13:04", "215C_EL PATO", 77, "LP1912")
    ,@(149, "11:34:59", "13:12", "16_SANTA ANA", 98, "LP1912")
    ,@(150, "11:47:17", "13:13", "16_SANTA ANA", 86, "LP1912")
    ,@(151, "11:34:59", "13:16", "10_OLMOS", 102, "LP1912")
    ,@(152, "11:47:17", "13:17", "10_OLMOS", 90, "LP1912")
    ,@(153, "11:54:18", "13:22", "23_HERNANDEZ", 88, "LP1912")
    ,@(154, "11:34:59", "13:24", "16_P MOR-SANTA ANA", 110, "LP1912")
    ,@(155, "11:47:17", "13:25", "16_P MOR-SANTA ANA", 98, "LP1912")
    ,@(156, "12:11:52", "13:25", "23_HERNANDEZ", 74, "LP1912")
    ,@(157, "11:34:59", "13:32", "215A_EL PATO", 118, "LP1912")
    ,@(158, "12:11:52", "13:32", "14_ABASTO", 81, "LP1912")
    ,@(159, "11:47:17", "13:33", "215A_EL PATO", 106, "LP1912")
    ,@(160, "12:11:52", "13:46", "225_GOMEZ", 95, "LP1912")
    ,@(161, "11:52:01", "13:47", "225_GOMEZ", 115, "LP1912")
    ,@(162, "12:11:52", "14:01", "10_OLMOS", 110, "LP1912")
)
Set-RowData $ws1 $rows1

# --- Sheet "LP1912-215" ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$rows2 = @(
    ,@(2, "Última actualización: 12:11:52", $null, $null, $null, $null)
)
Set-RowData $ws2 $rows2

# --- Sheet "6203-6173" ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$rows3 = @(
    ,@(2, "Última actualización: 12:11:52", $null, $null, $null, $null)
    ,@(3, "Total filas: 22", $null, $null, $null, $null)
    ,@(27, "12:11:52", "13:56", "215C_LA PLATA", 105, "L6203")
)
Set-RowData $ws3 $rows3

